# Add "comment" descriptions (column D) for each club row and fix a
# typo in the "Баштау" club name (-> "Бештау"), matching commit
# "Add 04 09 + mt @SvyTo".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the misspelled club name used in the answer-option columns.
$ws.Cells.Replace("Баштау", "Бештау")

# Russian-language descriptions ("comment" column, D) for each of the
# 19 club rows (rows 2-20), in row order.
$comments = @(
    "«Лада-Тольятти» — советский и российский футбольный клуб из города Тольятти, Самарская область. В 2010 и 2022 годах был расформирован, но в обоих случаях спустя два года возобновлял активную деятельность.",
    "«Жемчужина» — российский футбольный клуб из города Сочи. Выступал на профессиональном уровне с 1991 по 2003 и с 2008 по 2011 годы.",
    "«Тосно» — российский футбольный клуб из Тосно, Ленинградская область. Основан в 2013 году. Обладатель Кубка России сезона 2017/18. Расформирован в 2018 году. Возрождён в 2023.",
    "«Уралан» — российский футбольный клуб из Элисты. По окончании сезона 2004 года, по итогам которого клуб вылетел во Второй дивизион, отказался от участия в соревнованиях  и был лишён лицензии.",
    "«Лада» — российский футбольный клуб из города Димитровграда Ульяновской области, существовавший с 2017 до 2021 года. В сезонах 2019/20 и 2020/21 выступал в Первенстве ПФЛ.",
    "«Бештау» — бывший российский футбольный клуб из города Лермонтова, существовавший с 1992 по 2000 год. Лучшее достижение в первенстве России — 10 место в 1 зоне второй лиги в 1992 году.",
    "«Саранск» — бывший российский профессиональный футбольный клуб из одноимённого города, выступавший во Втором дивизионе ФНЛ в сезоне 2021/22, по завершении которого был расформирован.",
    "«Асмарал» — бывший советский и российский футбольный клуб из Москвы. Один из основателей Чемпионата России по футболу. В 2003 году клуб «Асмарал» был официально объявлен банкротом.",
    "«Текстильщик» — советский и российский футбольный клуб из города Камышина. С 1987 по 2008 год выступал в различных дивизионах чемпионата России. С 2009 участвует в первенстве Волгоградской области.",
    "«Арарат» — российский футбольный клуб из Москвы. Основан в 2017 году, расформирован 5 февраля 2020 года.",
    "«Автодор» — российский футбольный клуб из Владикавказа. В марте 2011 года ФК «Автодор» не смог пройти процедуру аттестации и снялся с первенства России по футболу среди команд второго дивизиона.",
    "«Кайрат» — российский футбольный клуб из Москвы, являвшийся фарм-клубом казахстанского «Кайрата». 11 июня 2022 года стало известно об объединении Кайрат-Москва с Кайрат-Жастар по финансовым причинам.",
    "«Кубань» — советский и российский футбольный клуб из Краснодара, существовавший с 1928 по 2018 год. На момент до расформирования был одним из старейших футбольных клубов России.",
    "«Луч» — бывший российский футбольный клуб из Владивостока. Основан в 1958 году. 1 апреля 2020 года было объявлено о прекращении выступлений в профессиональном футболе.",
    "«Ессентуки» — российский профессиональный футбольный клуб из одноимённого города Ставропольского края, основанный в 2016 году. С сезона 2020/21 выступал в группе 1 третьего по уровню дивизиона России.",
    "«Звезда́» — бывший советский и российский футбольный клуб из Перми. «Звезда» существовала с 1932 по 1996 года. Сезон 1996 года команда играла в первенстве Пермской области, после чего прекратила существование.",
    "«Волга» — бывший российский футбольный клуб из Нижнего Новгорода. Клуб был основан в 1998 году. 15 июня 2016 года руководство клуба приняло решение о расформировании клуба.",
    "«Тамбов» — российский футбольный клуб из Тамбова. В 2021 году Тамбов потерял шансы на сохранение прописки в РПЛ и руководство клуба объявило о том, что клуб лишится профессионального статуса",
    "«Звезда» — российский футбольный клуб из города Иркутска. Существовал в 1957—2008 годах."
)

$row = 2
foreach ($comment in $comments) {
    $ws.Cells.Item($row, 4).Value = $comment
    $row = $row + 1
}
